# Add a new "Compact List" paragraph style (styleId "CompactList"),
# mirroring the existing "Compact" style: based on Body Text, marked as
# a quick style, with 36-twip (1.8pt) spacing before/after.

$d = $word.ActiveDocument

$bodyText = $d.Styles("BodyText")

$compactList = $d.Styles.Add("CompactList", 1)
$compactList.NameLocal = "Compact List"
$compactList.BaseStyle = $bodyText
$compactList.QuickStyle = $true
$compactList.ParagraphFormat.SpaceBefore = 1.8
$compactList.ParagraphFormat.SpaceAfter = 1.8
